$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting the latest crypto price/volume scrape.
# Column D values are forced to Text so Excel does not reinterpret
# numeric-looking strings (e.g. "11.00", "0.0530") as numbers and
# strip precision / trailing zeros - matching the inline-string cells
# in the original workbook. Style is restored to "Normal" afterwards
# so no stray cell-style index is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.736.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.793.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.553"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.45"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.284"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0718"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0934"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.051.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.791.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.635"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.703.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0812"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.83%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0530"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("E31").Value = "  -2.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.437.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0192"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.12%  "
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.635"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "84.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.926"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.70%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.10%  "
$ws.Range("E46").Value = "  -4.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.950.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.06%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("E51").Value = "  +8.30%  "
